$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("July 2018")

# Mirror the existing "17-07-2018" rows (6 & 7) with three new FOCUS hearts
# 100 scale part requests added as rows 8, 9, 10.
$dateRequested = $ws.Range("A7").Value()
$material      = $ws.Range("E7").Value()
$shellCount    = $ws.Range("F7").Value()
$infill        = $ws.Range("G7").Value()
$layerHeight   = $ws.Range("H7").Value()
$comments      = $ws.Range("I7").Value()

$parts = @("FOCUS a4c 100 Scale", "FOCUS lax 100 Scale", "FOCUS sax 100 Scale")

$row = 8
foreach ($part in $parts) {
    $ws.Range("A$row").Value = $dateRequested
    $ws.Range("C$row").Value = $part
    $ws.Range("D$row").Value = 1
    $ws.Range("E$row").Value = $material
    $ws.Range("F$row").Value = $shellCount
    $ws.Range("G$row").Value = $infill
    $ws.Range("H$row").Value = $layerHeight
    $ws.Range("I$row").Value = $comments
    $row++
}

$ws.Range("E15").Select()
